$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Nov 1 2020 @ Abu Dhabi vs Chennai Super Kings" innings (formerly row 4)
# becomes the sole data row (row 2); all other match rows (old rows 2, 3, 5) go away.
$ws.Range("A2").Value = " Nov 1 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "Super Kings won by 9 wickets (with 7 balls remaining)"
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Mayank Agarwal "

# Numeric-looking stats are stored as text in this sheet, so force text
# entry (leading apostrophe, like typing '26 into Excel) then restore the
# Normal style so no stray formatting sticks to the cell.
$ws.Range("G2").Value = "'26"
$ws.Range("H2").Value = "'15"
$ws.Range("I2").Value = "'5"
$ws.Range("J2").Value = "'0"
$ws.Range("K2").Value = "'173.33"
$ws.Range("G2:K2").Style = "Normal"

# Drop the now-superseded rows 3-5 so only the header + the single match row remain.
$ws.Rows("3:5").Delete()
